$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("A1").Value = "Service Code"
$ws.Range("B1").Value = "Service Title"
$ws.Range("C1").Value = "Service Duration (minutes)"
$ws.Range("D1").ClearContents()

# --- Row 2 ---
$ws.Range("A2").Value = "D"
$ws.Range("B2").Value = "Counseling"
$ws.Range("C2").Value = 10
$ws.Range("D2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = "E"
$ws.Range("B3").Value = "Check-up"
$ws.Range("C3").Value = 8
$ws.Range("D3").ClearContents()

# --- Row 4 ---
$ws.Range("A4").Value = "F"
$ws.Range("B4").Value = "Surgery"
$ws.Range("C4").Value = 30
$ws.Range("D4").ClearContents()

# --- Rows 5-7: clear all data ---
$ws.Range("A5:D7").ClearContents()

# --- Re-apply the header style (s=1) to D1:D4 and to the now-empty rows 5:7 ---
$ws.Range("A1").Copy()
$ws.Range("D1:D4").PasteSpecial(-4122)
$ws.Range("A5:D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 1 height doubles (wrapped 2-line header) ---
$ws.Rows("1:1").RowHeight = 57.6

# --- Selection moves to A7 ---
$ws.Range("A7").Select() | Out-Null
